$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data.
# D-column values are plain text (e.g. "67.789.48"); force text format so
# Excel does not reinterpret them as numbers (avoiding float rounding) and
# reset the style afterwards so no extra cell formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.789.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.64%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.807.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.35%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.802.60"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  +0.17%  "

$ws.Range("E10").Value = "  +0.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.30"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.48%  "

$ws.Range("E12").Value = "  -0.59%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.446.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.820.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.782.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.85%  "

$ws.Range("E19").Value = "  +1.94%  "

$ws.Range("E20").Value = "  +0.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "461.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.89%  "

$ws.Range("E23").Value = "  +0.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000155"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.26%  "

$ws.Range("E26").Value = "  +1.82%  "

$ws.Range("E27").Value = "  -2.33%  "

$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("E29").Value = "  -0.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.953.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("E31").Value = "  +0.94%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.84%  "

$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("E36").Value = "  -0.32%  "

$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.48%  "

$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.74%  "

$ws.Range("E46").Value = "  -0.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.19%  "

$ws.Range("E48").Value = "  -0.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "392.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.86%  "

$ws.Range("E50").Value = "  -4.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.05%  "
